# Daily attendance processing - 2026-02-07 08:00:05 UTC
# Swap the "Last Name, First" style ordering of the Administrator's name in the
# "Recorded By" column (G): "Miss Dina Nasr, Administrator" -> "Administrator, Miss Dina Nasr"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "Miss Dina Nasr, Administrator"
$replacement = "Administrator, Miss Dina Nasr"

# Collect every matching cell address first (via Find/FindNext) so that we
# don't mutate the collection we are iterating over.
$addresses = New-Object System.Collections.Generic.List[string]

$first = $ws.Cells.Find($target)
if ($first -ne $null) {
    $firstAddr = $first.Address()
    [void]$addresses.Add($firstAddr)

    $current = $first
    do {
        $current = $ws.Cells.FindNext($current)
        if ($current -ne $null -and $current.Address() -ne $firstAddr) {
            [void]$addresses.Add($current.Address())
        }
    } while ($current -ne $null -and $current.Address() -ne $firstAddr)
}

foreach ($address in $addresses) {
    $ws.Range($address).Value = $replacement
}

Write-Host "Updated $($addresses.Count) cell(s) from '$target' to '$replacement'"
